$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of FX data for 2020-12-31
$ws.Range("A8").Value = 44196
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("B8").Value = "USD"
$ws.Range("C8").Value = "HKD"
$ws.Range("D8").Value = 7.7530999999999999

$ws.Range("D9").Select()
